$wb = $excel.ActiveWorkbook

$wsWeek5 = $wb.Worksheets.Item("Week 5")
$wsWeek6 = $wb.Worksheets.Item("Week 6")

# Make "Week 6" the active sheet/tab (was "Week 5"): updates workbook.xml
# activeTab and moves sheetView tabSelected from Week 5 to Week 6.
$wsWeek6.Activate()

# New row 2: 2018-02-12, 12:30 PM - 1:00 PM, "Meet with client to get
# product info", 0.5 hours. Copy the date/time number formats from an
# existing timesheet row so we reuse the workbook's existing styles
# instead of minting new ones.
$wsWeek5.Range("A2").Copy()
$wsWeek6.Range("A2").PasteSpecial(-4122) # xlPasteFormats
$wsWeek5.Range("B2").Copy()
$wsWeek6.Range("B2").PasteSpecial(-4122)
$wsWeek5.Range("C2").Copy()
$wsWeek6.Range("C2").PasteSpecial(-4122)

$wsWeek6.Rows.Item(2).RowHeight = 18

$wsWeek6.Range("A2").Value = 41681
$wsWeek6.Range("B2").Value = 0.52083333333333337
$wsWeek6.Range("C2").Value = 0.54166666666666663
$wsWeek6.Range("D2").Value = "Meet with client to get product info"
$wsWeek6.Range("E2").Value = 0.5

# New row 3: 2018-02-12, 4:30 PM - 5:30 PM, "Photographed and started
# editing product images", 1 hour.
$wsWeek5.Range("A2").Copy()
$wsWeek6.Range("A3").PasteSpecial(-4122)
$wsWeek5.Range("B2").Copy()
$wsWeek6.Range("B3").PasteSpecial(-4122)
$wsWeek5.Range("C2").Copy()
$wsWeek6.Range("C3").PasteSpecial(-4122)

$wsWeek6.Rows.Item(3).RowHeight = 18

$wsWeek6.Range("A3").Value = 41681
$wsWeek6.Range("B3").Value = 0.6875
$wsWeek6.Range("C3").Value = 0.72916666666666663
$wsWeek6.Range("D3").Value = "Photographed and started editing product images"
$wsWeek6.Range("E3").Value = 1

# Put the selection where the author left it on "Week 6" after entering
# the new rows.
$wsWeek6.Range("A4").Select()

$excel.CutCopyMode = $false
